# Reporte Consolidado Diario: fill in the "tipos de depositos" counts for
# the two inspector rows (13 & 14) and wire up the column totals in row 37
# (the "Totales" row) with SUM formulas over the data block (rows 13-37).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 13: FERNANDEZ MAURICIO LORENZO ----------------------------------
$ws.Range("A13").Value = 13

$cols = @('C','D','E','F','G','H','I','J','K','L','M','N','O','P','Q','R','S','T','U','V','W','X','Y','Z', `
          'AA','AB','AC','AD','AE','AF','AG','AH','AI','AJ','AK','AL','AM','AN','AO','AP','AQ','AR','AS','AT','AU','AV')

$row13 = @(14,1,0,0,0,0,0,0,0,0,0,5,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,5)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "13").Value = $row13[$i]
}

# --- Row 14: SDFSDFDS ------------------------------------------------------
$ws.Range("A14").Value = 14

$row14 = @(6,0,0,2,0,0,0,29,8,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,9)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "14").Value = $row14[$i]
}

# --- Row 37: column totals --------------------------------------------------
# "I" / "P" / "TQ" / "TF" / "D" depósito-type totals, columns J..AU, each
# summing its own column across the data rows (13..37).
$totalCols = @('J','K','L','M','N','O','P','Q','R','S','T','U','V','W','X','Y','Z', `
               'AA','AB','AC','AD','AE','AF','AG','AH','AI','AJ','AK','AL','AM','AN','AO','AP','AQ','AR','AS','AT','AU')

foreach ($col in $totalCols) {
    $ws.Range($col + "37").Formula = "=SUM(" + $col + "13:" + $col + "37)"
}
